$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5611346983557807
$ws.Range("C2").Value = -0.6242536099627518
$ws.Range("D2").Value = 0.5082040332081577

$ws.Range("B3").Value = -0.6761097870911426
$ws.Range("C3").Value = -0.6533155315349959
$ws.Range("D3").Value = -0.6302161934269627

$ws.Range("B4").Value = 0.7697467008299035
$ws.Range("C4").Value = -0.6195852097269362
$ws.Range("D4").Value = 0.6822729943094337

$ws.Range("B5").Value = -0.559148887524349
$ws.Range("C5").Value = 0.6724377503136852
$ws.Range("D5").Value = 0.6855731720921452

$ws.Range("B6").Value = -0.5711540419562121
$ws.Range("C6").Value = 0.5945815665727221
$ws.Range("D6").Value = -0.6857643802455191

$ws.Range("B7").Value = 0.7150708283385313
$ws.Range("C7").Value = -0.7073208358945109
$ws.Range("D7").Value = 0.7847531275613041

$ws.Range("B8").Value = -0.7061689258161673
$ws.Range("C8").Value = 0.8013238183578082
$ws.Range("D8").Value = -0.8102030208510096

$ws.Range("B9").Value = -0.6183026900085103
$ws.Range("C9").Value = 0.6171546141700381
$ws.Range("D9").Value = -0.6380200697984402
